$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.331.98'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '2.076.34'
$ws.Range("E3").Value = '  +4.43%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.60'
$ws.Range("E5").Value = '  -1.82%  '
$ws.Range("E6").Value = '  +2.18%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.34'
$ws.Range("E8").Value = '  +4.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.380'
$ws.Range("E9").Value = '  +2.33%  '
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0760'
$ws.Range("E11").Value = '  +1.44%  '
$ws.Range("E12").Value = '  +3.34%  '
$ws.Range("D13").Value = '2.381.60'
$ws.Range("E13").Value = '  +4.31%  '
$ws.Range("E14").Value = '  +3.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.02'
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.775'
$ws.Range("E16").Value = '  +2.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.24'
$ws.Range("E17").Value = '  +3.92%  '
$ws.Range("D18").Value = '2.072.94'
$ws.Range("E18").Value = '  +4.06%  '
$ws.Range("D19").Value = '37.448.04'
$ws.Range("E19").Value = '  +1.36%  '
$ws.Range("E20").Value = '  +20.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.40'
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '223.85'
$ws.Range("E23").Value = '  -1.85%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  +3.04%  '
$ws.Range("E26").Value = '  +1.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.85'
$ws.Range("E27").Value = '  +1.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.87'
$ws.Range("E28").Value = '  +2.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.129'
$ws.Range("E29").Value = '  +6.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.31'
$ws.Range("E30").Value = '  +1.22%  '
$ws.Range("E31").Value = '  +7.40%  '
$ws.Range("E32").Value = '  +1.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.45'
$ws.Range("E33").Value = '  +1.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0622'
$ws.Range("E34").Value = '  +1.56%  '
$ws.Range("E35").Value = '  +9.28%  '
$ws.Range("E36").Value = '  +4.45%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.96'
$ws.Range("E38").Value = '  +14.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.33'
$ws.Range("E39").Value = '  +0.75%  '
$ws.Range("E40").Value = '  -1.31%  '
$ws.Range("E41").Value = '  -2.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.48'
$ws.Range("E42").Value = '  +22.90%  '
$ws.Range("E43").Value = '  +8.85%  '
$ws.Range("D44").Value = '1.471.51'
$ws.Range("E44").Value = '  +2.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '95.05'
$ws.Range("E45").Value = '  +7.68%  '
$ws.Range("E46").Value = '  +3.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.11'
$ws.Range("E47").Value = '  +5.59%  '
$ws.Range("E48").Value = '  +1.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.03'
$ws.Range("E49").Value = '  +2.71%  '
$ws.Range("E50").Value = '  +7.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.94'
$ws.Range("E51").Value = '  +2.21%  '
